# Applies the "Describe..." -> "Contain..." rewrites in the SRS table
# (ERD / Entity Data Dictionary section) as described by the commit.

$d = $word.ActiveDocument

$apos = [char]0x2019

$pairs = @(
    @("Describe the customer${apos}s server in data center.", "Contain the server information in data center."),
    @("Describe all of statuses of objects in data center.", "Contain all of statuses of objects in data center."),
    @("Describe all location in data center.", "Contain all location information in data center."),
    @("Describe all racks which are putting in data center.", "Contain all racks information which are putting in data center."),
    @("Describe all logs about object${apos}s changes in data center.", "Contain all logs about object${apos}s changes in data center."),
    @("Describe all roles in the system.", "Contain all roles in the system."),
    @("Describe all of current IP Addresses of server.", "Contain all of current IP Addresses of server."),
    @("Describe content of each request which was sent by customer.", "Contain all request content which was sent by customer."),
    @("Describe all user${apos}s accounts in the system.", "Contain all user information in the system."),
    @("Describe all note which was wrote by previous shift for the next shift.", "Contain all note which was wrote by previous shift for the next shift."),
    @("Describe temporary detail of all requests.", "Contain temporary detail of all requests."),
    @("Describe all IP Addresses which data center is keeping.", "Contain all IP Addresses information which data center is keeping."),
    @("Describe all types of request.", "Contain all types of request."),
    @("Describe all types of log.", "Contain all types of log."),
    @("Describe all shift group of data center. ", "Contain all shift group of data center. "),
    @("Describe which group is in which shift each day.", "Contain which group is in which shift each day."),
    @("Describe started time and ended time of each shift group.", "Contain started time and ended time of each shift group."),
    @("Describe all racks which was rent by customer.", "Contain all racks which was rent by customer."),
    @("Describe all contents of each task.", "Contain all contents of each task.")
)

foreach ($pair in $pairs) {
    $find = $pair[0]
    $replace = $pair[1]
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# This cell's sentence-ending period sits in its own run *after* a
# "_GoBack" bookmark, so the replacement is deliberately done without the
# trailing period to avoid swallowing (and deleting) that bookmark.
$d.Content.Find.Execute("Describe all contents of each notification", $true, $false, $false, $false, $false, $true, 1, $false, "Contain all notification information", 2)
